$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows (row, A..O) ---------------------------------------------
# Column order: A, B, C(text), D, E, F, G, H, I, J, K, L, M, N, O
$rows = @(
  @(1754, 39323, "CEI- SubidoVie10_11", 1, 87654325, 1, 110, 1, 20191103, 84043, 0, 401,   71,   0, 0),
  @(1755, 39323, "CEI- SubidoVie10_12", 1, 87654325, 1, 110, 1, 20191203, 84043, 0, 401,   71,   0, 0),
  @(1756, 39323, "CEI- SubidoVie10_12", 1, 87654323, 1, 110, 1, 20191203, 84043, 0, 401,   71,   0, 0),
  @(1757, 39323, "CEI- SubidoVie10_01", 1, 87654325, 1, 110, 1, 20200103, 84043, 0, 401,   71,   0, 0),
  @(1758, 39323, "CEI- SubidoVie10_01", 1, 87654323, 1, 110, 1, 20200103, 84043, 0, 401,   71,   0, 0),
  @(1759, 39323, "CEI- SubidoVie10_01", 1, 87654321, 1, 110, 1, 20200103, 84043, 0, 401,   71,   0, 0),
  @(1760, 39323, "CEI- SubidoVie10_01", 1, 87654324, 1, 110, 1, 20200103, 84043, 0, 401,   71,   0, 0),
  @(1981, 39323, "CEI- SubidoVie10_02", 1, 87654325, 1, 110, 0, 20200207, 122113, 0, 74, 2186, 0, 0),
  @(2655, 39323, "CEI- SubidoVie10_02", 1, 87654323, 1, 110, 0, 20200209, 122113, 0, 74, 2186, 0, 0),
  @(3719, 39323, "CEI- SubidoVie10_02", 1, 87654321, 1, 110, 0, 20200204, 83224,  0, 24,  340,  0, 0),
  @(4338, 39323, "CEI- SubidoVie10_02", 1, 87654324, 1, 110, 0, 20200203, 84043,  0, 401,  71,  0, 0)
)

# --- Build the formatting (wrapText + vertical-center) on A2 first, as a
# single clean style, then copy/paste-special that format onto the full
# data range so every row shares the same one new cellXfs entry. ---------
$firstCell = $ws.Cells.Item(2, 1)
$firstCell.WrapText = $true
$firstCell.VerticalAlignment = -4108
$firstCell.Copy()
$dataRange = $ws.Range("A2:O12")
$dataRange.PasteSpecial(-4122)

# --- Write the values/formulas for every data row -----------------------
$r = 2
foreach ($row in $rows) {
  for ($col = 1; $col -le 15; $col++) {
    $ws.Cells.Item($r, $col).Value = $row[$col - 1]
  }
  $r++
}

# --- Selection, matching the recorded end-state --------------------------
[void]$ws.Range("G17").Select()
